{"js": "// Update the worksheet date heading and the 25 division-problem cells\n// in the single practice table. Every value is addressed by its fixed\n// (row, column) position in the table (not by its old text) because a\n// handful of the original expressions are duplicated across cells but\n// map to different replacement text, so a plain global find/replace\n// would be ambiguous/unsafe.\n\nconst body = context.document.body;\n\n// --- 1. Heading paragraph: \"2025-08-01 Friday\" -> \"2025-08-02 Saturday\"\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst heading = paragraphs.items[0];\nheading.insertText(\"2025-08-02 Saturday\", Word.InsertLocation.replace);\n\n// --- 2. Table cells: replace the division problem in each populated cell.\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// [rowIndex, columnIndex, newText]\nconst updates = [\n  [0, 0, \"64\u00f76=10, 4\"],\n  [0, 1, \"48\u00f74=12, 0\"],\n  [0, 2, \"66\u00f76=11, 0\"],\n  [0, 3, \"52\u00f79=5, 7\"],\n  [0, 4, \"41\u00f77=5, 6\"],\n\n  [4, 0, \"28\u00f73=9, 1\"],\n  [4, 1, \"37\u00f72=18, 1\"],\n  [4, 2, \"37\u00f75=7, 2\"],\n  [4, 3, \"29\u00f79=3, 2\"],\n  [4, 4, \"36\u00f78=4, 4\"],\n\n  [8, 0, \"77\u00f75=15, 2\"],\n  [8, 1, \"78\u00f73=26, 0\"],\n  [8, 2, \"85\u00f74=21, 1\"],\n  [8, 3, \"65\u00f75=13, 0\"],\n  [8, 4, \"35\u00f76=5, 5\"],\n\n  [12, 0, \"56\u00f74=14, 0\"],\n  [12, 1, \"91\u00f75=18, 1\"],\n  [12, 2, \"65\u00f78=8, 1\"],\n  [12, 3, \"49\u00f77=7, 0\"],\n  [12, 4, \"54\u00f74=13, 2\"],\n\n  [16, 0, \"63\u00f76=10, 3\"],\n  [16, 1, \"37\u00f72=18, 1\"],\n  [16, 2, \"31\u00f73=10, 1\"],\n  [16, 3, \"38\u00f73=12, 2\"],\n  [16, 4, \"82\u00f78=10, 2\"],\n];\n\nfor (const [rowIndex, colIndex, newText] of updates) {\n  const cell = table.getCell(rowIndex, colIndex);\n  const firstPara = cell.body.paragraphs.getFirst();\n  firstPara.insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Update the worksheet date heading and the 25 division-problem cells\n# in the single practice table. Every value is addressed by its fixed\n# (row, column) position in the table (not by its old text) because a\n# handful of the original expressions are duplicated across cells but\n# map to different replacement text, so a plain global find/replace\n# would be ambiguous/unsafe.\n#\n# Word COM's Table.Cell(row, column) is 1-based.\n\n$d = $word.ActiveDocument\n\n# --- 1. Heading paragraph: \"2025-08-01 Friday\" -> \"2025-08-02 Saturday\"\n$heading = $d.Paragraphs.Item(1)\n$heading.Range.Text = \"2025-08-02 Saturday\"\n\n# --- 2. Table cells: replace the division problem in each populated cell.\n$t = $d.Tables.Item(1)\n\n$updates = @(\n    @(1, 1, \"64\u00f76=10, 4\"),\n    @(1, 2, \"48\u00f74=12, 0\"),\n    @(1, 3, \"66\u00f76=11, 0\"),\n    @(1, 4, \"52\u00f79=5, 7\"),\n    @(1, 5, \"41\u00f77=5, 6\"),\n\n    @(5, 1, \"28\u00f73=9, 1\"),\n    @(5, 2, \"37\u00f72=18, 1\"),\n    @(5, 3, \"37\u00f75=7, 2\"),\n    @(5, 4, \"29\u00f79=3, 2\"),\n    @(5, 5, \"36\u00f78=4, 4\"),\n\n    @(9, 1, \"77\u00f75=15, 2\"),\n    @(9, 2, \"78\u00f73=26, 0\"),\n    @(9, 3, \"85\u00f74=21, 1\"),\n    @(9, 4, \"65\u00f75=13, 0\"),\n    @(9, 5, \"35\u00f76=5, 5\"),\n\n    @(13, 1, \"56\u00f74=14, 0\"),\n    @(13, 2, \"91\u00f75=18, 1\"),\n    @(13, 3, \"65\u00f78=8, 1\"),\n    @(13, 4, \"49\u00f77=7, 0\"),\n    @(13, 5, \"54\u00f74=13, 2\"),\n\n    @(17, 1, \"63\u00f76=10, 3\"),\n    @(17, 2, \"37\u00f72=18, 1\"),\n    @(17, 3, \"31\u00f73=10, 1\"),\n    @(17, 4, \"38\u00f73=12, 2\"),\n    @(17, 5, \"82\u00f78=10, 2\")\n)\n\nforeach ($u in $updates) {\n    $row = $u[0]\n    $col = $u[1]\n    $text = $u[2]\n    $cell = $t.Cell($row, $col)\n    $cell.Range.Text = $text\n}\n"}
